# Updates to address #5
# EvennessWindReductions.xlsx - update RangeMaximum values and move the
# active selection down one row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# RangeMaximum column: second bucket boundary grows from 0.5 to 25,
# third bucket boundary grows from 1 to 50.
$ws.Range("A2").Value = 25
$ws.Range("A3").Value = 50

# The cursor ends up resting on A4 after editing A3.
$ws.Activate()
$ws.Range("A4").Select()
